# SequentialExtractions.xlsx edit: append new t1/t0/t2/t3/ts0 BD+NaOH rows to
# the "sediment" sheet, fix a doubled D14 value, bold the relevant label /
# value columns for the freshly-entered rows, and re-point the active
# sheet/selection at sheet1.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "sediment"
$ws2 = $wb.Worksheets.Item(2)   # "ferrosorb"

# --- fix the mis-entered (halved) value in D14 -----------------------------
$ws1.Cells.Item(14, 4).Value = 18.784662925916898

# --- append the new measurement rows (50-79) -------------------------------
$newRows = @(
    @("t1",  "BD",   1, 12.5,               0,    0.4783,              "Fe"),
    @("t1",  "BD",   2, 12.1,               0,    0.47499999999999998, "Fe"),
    @("t1",  "BD",   3, 15.7,               0,    0.62560000000000004, "Fe"),
    @("t0",  "BD",   1, 0.89,               0,    0.43219999999999997, "Fe"),
    @("t0",  "BD",   2, 1.08,               0,    0.58650000000000002, "Fe"),
    @("t0",  "BD",   3, 0.99,               0,    0.53480000000000005, "Fe"),
    @("t2",  "BD",   1, 6.98,               0,    0.53100000000000003, "Fe"),
    @("t2",  "BD",   2, 7.25,               0,    0.52929999999999999, "Fe"),
    @("t2",  "BD",   3, 7.33,               0,    0.51039999999999996, "Fe"),
    @("t3",  "BD",   1, 8.16,               0,    0.58209999999999995, "Fe"),
    @("t3",  "BD",   2, 8.34,               0,    0.57420000000000004, "Fe"),
    @("t3",  "BD",   3, 7.89,               0,    0.53310000000000002, "Fe"),
    @("t1",  "NaOH", 1, 0.79,               0.22, 0.4783,              "Fe"),
    @("t1",  "NaOH", 2, 0.71,               0.22, 0.47499999999999998, "Fe"),
    @("t1",  "NaOH", 3, 0.85,               0.22, 0.62560000000000004, "Fe"),
    @("t0",  "NaOH", 1, 0.6,                0.22, 0.43219999999999997, "Fe"),
    @("t0",  "NaOH", 2, 0.68,               0.22, 0.58650000000000002, "Fe"),
    @("t0",  "NaOH", 3, 0.64,               0.22, 0.53480000000000005, "Fe"),
    @("t2",  "NaOH", 1, 1.1000000000000001, 0.21, 0.53100000000000003, "Fe"),
    @("t2",  "NaOH", 2, 1.25,               0.21, 0.52929999999999999, "Fe"),
    @("t2",  "NaOH", 3, 1.22,               0.21, 0.51039999999999996, "Fe"),
    @("t3",  "NaOH", 1, 1.39,               0.21, 0.58209999999999995, "Fe"),
    @("t3",  "NaOH", 2, 1.34,               0.21, 0.57420000000000004, "Fe"),
    @("t3",  "NaOH", 3, 1.31,               0.21, 0.53310000000000002, "Fe"),
    @("ts0", "BD",   1, 6.44,               0,    0.48948000000000003, "Fe"),
    @("ts0", "BD",   2, 5.44,               0,    0.41206999999999999, "Fe"),
    @("ts0", "BD",   3, 5.35,               0,    0.42157,             "Fe"),
    @("ts0", "NaOH", 1, 1.36,               0.21, 0.48948000000000003, "Fe"),
    @("ts0", "NaOH", 2, 1.23,               0.21, 0.41206999999999999, "Fe"),
    @("ts0", "NaOH", 3, 1.25,               0.21, 0.42157,             "Fe")
)

$startRow = 50
$r = $startRow
foreach ($row in $newRows) {
    for ($c = 0; $c -lt 7; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $row[$c]
    }
    $r++
}
$endRow = $r - 1   # 79

# --- bold formatting for the new block --------------------------------------
# Column A (Fraction label) is bold for every new row except the "ts0" rows.
$ws1.Range("A50:A73").Font.Bold = $true

# For the "BD" sub-block (rows 50-61) the value/parameter columns are bold too.
$ws1.Range("D50:E61").Font.Bold = $true
$ws1.Range("G50:G61").Font.Bold = $true

# --- re-point active sheet / selection --------------------------------------
# (sheet2's own selection, O17, is unchanged - only the active sheet/tab and
# sheet1's selected cell move.)
$ws1.Activate()
$ws1.Range("I16").Select()
